# Add a new "2022" column (Q) to the yearly table, mirroring the
# formatting already used by the adjacent "2021" column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: year label 2022, formatted like P4 (2021).
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2022

# Data cell: value 64.2, formatted like P5 (80.9).
$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 64.2

# Match the saved selection/active cell from the source workbook.
$ws.Range("R4").Select()
